$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as plain text so values (e.g. "65.30") are not
# auto-coerced into numbers and lose formatting/trailing zeros.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.697.66'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.599.61'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.34'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0843'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.824.11'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.620.80'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.04'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.30'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.678.03'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0756'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '209.60'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.17'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.87%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.94'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.18'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.49%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.56%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.289.94'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.619'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.32%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.91%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.06'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +15.92%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.11'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.736.81'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.97'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.37'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.98%  '
